# Update "想去人数" (column F) counts across the workbook sheets to reflect
# a fresh scrape of the source data, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 469
$ws1.Range("F6").Value  = 70
$ws1.Range("F8").Value  = 350
$ws1.Range("F9").Value  = 1761
$ws1.Range("F10").Value = 379
$ws1.Range("F11").Value = 1435
$ws1.Range("F12").Value = 821
$ws1.Range("F13").Value = 346
$ws1.Range("F14").Value = 691
$ws1.Range("F15").Value = 12902
$ws1.Range("F16").Value = 12870
$ws1.Range("F19").Value = 11
$ws1.Range("F20").Value = 526
$ws1.Range("F21").Value = 54
$ws1.Range("F22").Value = 590
$ws1.Range("F25").Value = 14
$ws1.Range("F26").Value = 5
$ws1.Range("F28").Value = 98
$ws1.Range("F30").Value = 687

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 83

# Sheet "本地生活" (sheetId 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 173

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 173
$ws4.Range("F6").Value  = 469
$ws4.Range("F9").Value  = 70
$ws4.Range("F13").Value = 350
$ws4.Range("F14").Value = 1761
$ws4.Range("F15").Value = 379
$ws4.Range("F16").Value = 1435
$ws4.Range("F17").Value = 821
$ws4.Range("F18").Value = 346
$ws4.Range("F20").Value = 691
$ws4.Range("F21").Value = 12902
$ws4.Range("F22").Value = 12870
$ws4.Range("F23").Value = 963
$ws4.Range("F25").Value = 11
$ws4.Range("F26").Value = 526
$ws4.Range("F27").Value = 54
$ws4.Range("F28").Value = 590
$ws4.Range("F33").Value = 14
$ws4.Range("F34").Value = 5
$ws4.Range("F38").Value = 98
$ws4.Range("F40").Value = 687
$ws4.Range("F41").Value = 83
